# Insert a new data row at row 509 on the active sheet (a weekly "Papa"
# price record for Feria Lagunitas de Puerto Montt), shifting the
# existing rows 509:568 down to 510:569.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 509 and below down by one row.
$ws.Rows.Item(509).Insert()

# Populate the newly inserted row with the new record.
$row = 509
$ws.Cells.Item($row, 1).Value  = 4
$ws.Cells.Item($row, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value  = "Los Lagos"
$ws.Cells.Item($row, 4).Value2 = 44918
$ws.Cells.Item($row, 5).Value  = 10
$ws.Cells.Item($row, 6).Value  = 100114001
$ws.Cells.Item($row, 7).Value  = "Papa"
$ws.Cells.Item($row, 8).Value  = "Patagonia"
$ws.Cells.Item($row, 9).Value  = "1a nueva(o)"
$ws.Cells.Item($row, 10).Value = 600
$ws.Cells.Item($row, 11).Value = 16000
$ws.Cells.Item($row, 12).Value = 16000
$ws.Cells.Item($row, 13).Value = 16000
$ws.Cells.Item($row, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item($row, 15).Value = "Región de La Araucanía"
$ws.Cells.Item($row, 16).Value = 640
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
